$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H92").Value = 423.44446
$ws.Range("I92").Value = 404.95456
$ws.Range("K92").Value = 404.95456
$ws.Range("M92").Value = 843.04544
$ws.Range("H98").Value = 3600.5386
$ws.Range("I98").Value = 3775.5833
$ws.Range("K98").Value = 3775.5833
$ws.Range("M98").Value = -2277.5833
$ws.Range("H106").Value = 3567.3333
$ws.Range("I106").Value = 3567.3333
$ws.Range("K106").Value = 3567.3333
$ws.Range("M106").Value = -2936.3333
$ws.Range("H122").Value = 3600.5386
$ws.Range("I122").Value = 3775.5833
$ws.Range("K122").Value = 11326.7499
$ws.Range("M122").Value = -8876.749899999999
$ws.Range("H135").Value = 4258.5
$ws.Range("I135").Value = 3999.6667
$ws.Range("J135").Value = 5035
$ws.Range("K135").Value = 35997.0003
$ws.Range("L135").Value = 45315
$ws.Range("M135").Value = -33462.0003
$ws.Range("N135").Value = -50385
$ws.Range("H141").Value = 4173.4
$ws.Range("I141").Value = 2217.1667
$ws.Range("J141").Value = 11998.333
$ws.Range("K141").Value = 6651.500100000001
$ws.Range("L141").Value = 35994.999
$ws.Range("M141").Value = -1471.500100000001
$ws.Range("N141").Value = -46354.999

# --- ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 5615.125
$ws.Range("I32").Value = 5225.9546
$ws.Range("J32").Value = 6471.3
$ws.Range("K32").Value = 5225.9546
$ws.Range("L32").Value = 6471.3
$ws.Range("M32").Value = -4938.9546
$ws.Range("N32").Value = -7045.3
$ws.Range("H61").Value = 2822
$ws.Range("I61").Value = 1842.2106
$ws.Range("K61").Value = 1842.2106
$ws.Range("M61").Value = -1630.2106
$ws.Range("H63").Value = 2560.4285
$ws.Range("I63").Value = 2706.3333
$ws.Range("J63").Value = 2297.8
$ws.Range("K63").Value = 2706.3333
$ws.Range("L63").Value = 2297.8
$ws.Range("M63").Value = -2020.3333
$ws.Range("N63").Value = -3669.8
$ws.Range("H66").Value = 2560.4285
$ws.Range("I66").Value = 2706.3333
$ws.Range("J66").Value = 2297.8
$ws.Range("K66").Value = 13531.6665
$ws.Range("L66").Value = 11489
$ws.Range("M66").Value = -10099.6665
$ws.Range("N66").Value = -18353
$ws.Range("H74").Value = 106966.36
$ws.Range("I74").Value = 180276.39
$ws.Range("J74").Value = 3665.8635
$ws.Range("K74").Value = 180276.39
$ws.Range("L74").Value = 3665.8635
$ws.Range("M74").Value = -179402.39
$ws.Range("N74").Value = -5413.863499999999
$ws.Range("H77").Value = 106966.36
$ws.Range("I77").Value = 180276.39
$ws.Range("J77").Value = 3665.8635
$ws.Range("K77").Value = 901381.9500000001
$ws.Range("L77").Value = 18329.3175
$ws.Range("M77").Value = -897013.9500000001
$ws.Range("N77").Value = -27065.3175
$ws.Range("H97").Value = 2136.625
$ws.Range("I97").Value = 1975.9231
$ws.Range("K97").Value = 1975.9231
$ws.Range("M97").Value = -1479.9231
$ws.Range("H102").Value = 5545.6924
$ws.Range("I102").Value = 5819.5
$ws.Range("K102").Value = 5819.5
$ws.Range("M102").Value = -4197.5
$ws.Range("H122").Value = 2695
$ws.Range("I122").Value = 2679.75
$ws.Range("K122").Value = 8039.25
$ws.Range("M122").Value = -5589.25
$ws.Range("H132").Value = 2840.1785
$ws.Range("I132").Value = 1953.5714
$ws.Range("K132").Value = 5860.7142
$ws.Range("M132").Value = -3330.7142
$ws.Range("H136").Value = 2822
$ws.Range("I136").Value = 1842.2106
$ws.Range("K136").Value = 5526.6318
$ws.Range("M136").Value = -2976.6318

# --- BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H103").Value = 19576.25
$ws.Range("J103").Value = 19576.25
$ws.Range("L103").Value = 19576.25
$ws.Range("N103").Value = -21920.25

# --- CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H94").Value = 2252.7144
$ws.Range("J94").Value = 2443.5557
$ws.Range("L94").Value = 2443.5557
$ws.Range("N94").Value = -3345.5557
$ws.Range("H132").Value = 2308.1904
$ws.Range("I132").Value = 1641.2858
$ws.Range("K132").Value = 4923.857400000001
$ws.Range("M132").Value = -2393.857400000001

# --- CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 1812872.8
$ws.Range("I4").Value = 1185316.4
$ws.Range("K4").Value = 3555949.2
$ws.Range("M4").Value = -3555837.2
$ws.Range("H8").Value = 886.3333
$ws.Range("I8").Value = 886.3333
$ws.Range("K8").Value = 2658.9999
$ws.Range("M8").Value = -2519.9999
$ws.Range("H12").Value = 176.18182
$ws.Range("J12").Value = 176.18182
$ws.Range("L12").Value = 528.5454599999999
$ws.Range("N12").Value = -874.5454599999999
$ws.Range("H131").Value = 2114.9443
$ws.Range("J131").Value = 2003.6666
$ws.Range("L131").Value = 6010.9998
$ws.Range("N131").Value = -16090.9998
$ws.Range("H136").Value = 897.5
$ws.Range("I136").Value = 897.5
$ws.Range("K136").Value = 2692.5
$ws.Range("M136").Value = 2407.5
$ws.Range("H137").Value = 3096.8823
$ws.Range("I137").Value = 3072
$ws.Range("K137").Value = 9216
$ws.Range("M137").Value = -4116

# --- GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 50002016
$ws.Range("I80").Value = 90910824
$ws.Range("J80").Value = 2362.111
$ws.Range("K80").Value = 90910824
$ws.Range("L80").Value = 2362.111
$ws.Range("M80").Value = -90909826
$ws.Range("N80").Value = -4358.111
$ws.Range("H83").Value = 50002016
$ws.Range("I83").Value = 90910824
$ws.Range("J83").Value = 2362.111
$ws.Range("K83").Value = 454554120
$ws.Range("L83").Value = 11810.555
$ws.Range("M83").Value = -454549128
$ws.Range("N83").Value = -21794.555
$ws.Range("H102").Value = 1912.3529
$ws.Range("I102").Value = 1899.3077
$ws.Range("J102").Value = 1954.75
$ws.Range("K102").Value = 1899.3077
$ws.Range("L102").Value = 1954.75
$ws.Range("M102").Value = -277.3077000000001
$ws.Range("N102").Value = -5198.75
$ws.Range("H107").Value = 878.8
$ws.Range("I107").Value = 800.3333
$ws.Range("K107").Value = 800.3333
$ws.Range("M107").Value = 1119.6667
$ws.Range("H130").Value = 23500
$ws.Range("J130").Value = 23500
$ws.Range("L130").Value = 23500
$ws.Range("N130").Value = -33540
$ws.Range("H132").Value = 3594.6086
$ws.Range("I132").Value = 3401.9285
$ws.Range("J132").Value = 3894.3333
$ws.Range("K132").Value = 10205.7855
$ws.Range("L132").Value = 11682.9999
$ws.Range("M132").Value = -7675.7855
$ws.Range("N132").Value = -16742.9999

# --- LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 6236.8057
$ws.Range("I40").Value = 5996.6665
$ws.Range("J40").Value = 6717.0835
$ws.Range("K40").Value = 5996.6665
$ws.Range("L40").Value = 6717.0835
$ws.Range("M40").Value = -5860.6665
$ws.Range("N40").Value = -6989.0835
$ws.Range("H122").Value = 3332.875
$ws.Range("I122").Value = 3196.75
$ws.Range("K122").Value = 9590.25
$ws.Range("M122").Value = -7140.25
$ws.Range("H132").Value = 6787.0557
$ws.Range("I132").Value = 2852.8
$ws.Range("K132").Value = 8558.400000000001
$ws.Range("M132").Value = -6028.400000000001

# --- WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 1592.7778
$ws.Range("I107").Value = 1226.4286
$ws.Range("K107").Value = 3679.2858
$ws.Range("M107").Value = -1759.2858
$ws.Range("H122").Value = 25001524
$ws.Range("I122").Value = 1580.25
$ws.Range("J122").Value = 125001304
$ws.Range("K122").Value = 4740.75
$ws.Range("L122").Value = 375003912
$ws.Range("M122").Value = -2290.75
$ws.Range("N122").Value = -375008812
